# Added Week 15 simulations
# Update Row 3 ("R" - Road) target depth data on both OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 296
$wsOff.Range("C3").Value = 190
$wsOff.Range("D3").Value = 138
$wsOff.Range("E3").Value = 54
$wsOff.Range("F3").Value = 4
$wsOff.Range("G3").Value = 10

# DEF sheet
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 424
$wsDef.Range("C3").Value = 325
$wsDef.Range("D3").Value = 104
$wsDef.Range("E3").Value = 51
$wsDef.Range("F3").Value = 4
